$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 1.136768169855944
$ws.Range("C8").Value = 58.13946918268717
$ws.Range("D8").Value = 4687.969884107692
$ws.Range("E8").Value = 400134.7069724926
$ws.Range("F8").Value = 31383316.41009646
$ws.Range("G8").Value = 1265644238.911851
$ws.Range("H8").Value = 62460288396.99023
$ws.Range("I8").Value = 1436397544180646
